# fix done login admin user
#
# Adds a new team member row (Le Thanh Dat / Inter / email / phone) to the
# bottom (previously-blank) row 19 of the contact sheet, wires up the
# mailto: hyperlink on the e-mail cell exactly like the other rows in the
# table, and moves the active selection/scroll position to reflect where
# the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data (row 19) ------------------------------------------------
$ws.Range("A19").Value = "Lê Thành Đạt "
$ws.Range("B19").Value = "Inter"
$ws.Range("C19").Value = "dat.lethanh2@ncc.asia"
$ws.Range("D19").Value = "01 291"

# --- Hyperlink the new e-mail address, matching the other rows --------------
$null = $ws.Hyperlinks.Add(
    $ws.Range("C19"),
    "mailto:dat.lethanh2@ncc.asia",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "dat.lethanh2@ncc.asia"
)

# --- Restore the view/selection state ----------------------------------------
$null = $ws.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 3
